# Add a "Save" column (H) to the s_vals sheet, matching the header
# formatting of the existing stat columns (B1:G1) and a numeric value
# of 1 for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "sum" header's formatting (bold font, borders, centered
# alignment) onto the new header cell, then overwrite its text/value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell: numeric value of 1 (no special formatting, matching
# the plain numeric cells B2:G2).
$ws.Range("H2").Value = 1
